# Meet_the_Family bug codex.xlsx -- integration & unit tests added
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new bold "section header" row above the first test row (old row 2),
#    pushing the existing MTF_UT_0001..MTF_UT_0018 block down by one row.
$ws.Rows.Item(2).Insert()

# Row 2 becomes the new "test_members" section header (bold, column A only).
$ws.Range("A2").Value = "test_members"
$ws.Range("A2:C2").Font.Bold = $true

# 2. Append two more integration-test rows after the existing block (which now
#    ends at row 20, since everything shifted down by one).
$ws.Range("A21").Value = "MTF_IT_0001"
$ws.Range("B21").Value = "test_set_methods"

$ws.Range("A22").Value = "MTF_IT_0002"
$ws.Range("B22").Value = "test_get_relationship_methods"

# 3. New bold section header row for the family-tree tests (column A only).
$ws.Range("A23").Value = "test_family_tree"
$ws.Range("A23").Font.Bold = $true

# 4. Final data row for the new family-tree unit test.
$ws.Range("A24").Value = "MTF_UT_0019"
$ws.Range("B24").Value = "test_initialization"

# Leave the cursor where the author last left it.
[void]$ws.Range("K27").Select()
